$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (old B -> D, old C -> E)
$ws.Columns("B:C").Insert()

# Keep the same (custom) column width of 8 characters on the columns that
# now hold the date-rating data (C, D, E) just like the original column C.
$ws.Columns("C:E").ColumnWidth = 7.1667

# New header row values (write C1 "Jun_15" first so that new shared
# strings are appended in the same order seen upstream: Jun_15, then Jun_17)
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the new B and C columns (rows 2-27) with the placeholder "UN" text,
# matching the rest of the sheet's "no rating change" placeholder.
$ws.Range("B2:C27").Value = "UN"

# Row 21 (Wells Fargo & Co) keeps its special rating-change text in the
# newly inserted column C, matching the value carried into column D.
$ws.Range("C21").Value = "6/13/2018,Reiterates,Hold,"
